$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Plan": the implementation checklist is extended with the new
# pseudo-code steps that cover computing the learner average.
# ---------------------------------------------------------------------
$plan = $wb.Worksheets.Item("Plan")

# Bold section headers (match style already used by B4/B12/B22).
$plan.Range("B4").Font.Bold = $true
$plan.Range("B12").Font.Bold = $true
$plan.Range("B22").Font.Bold = $true
$plan.Range("B37").Font.Bold = $true

# Replace the old single note under the "Calculate average" heading with
# the expanded set of implementation notes.
$plan.Range("C39").Value = "Array with Submission ID : Max Points"
$plan.Range("D39").Value = "(created when fetching the submission info and stored in an array to be rehused) "

$plan.Range("A40:A48").EntireRow.Insert()
$plan.Range("C40").Value = "Array with results (learners obj)"
$plan.Range("C41").Value = "Variable to store sum of score"
$plan.Range("C42").Value = "Variable to store sum of Max points"
$plan.Range("C43").Value = "Variable to store avg result "
$plan.Range("C45").Value = "take each student from the array Results and on each submissionID"
$plan.Range("C46").Value = "sum of score"
$plan.Range("C47").Value = "sum max points"
$plan.Range("C49").Value = "get avg and store in avg result"
$plan.Range("C51").Value = "Update leaner object with ""avg "", avg result"

# ---------------------------------------------------------------------
# Sheet "Data": tidy up the now-unused border-only helper styles on the
# empty D:H cells beneath the "LearnerSubmission" block (rows 28-36).
# ---------------------------------------------------------------------
$data = $wb.Worksheets.Item("Data")
$data.Range("D28:H36").ClearFormats()
